# Weekly refresh: reshuffle the per-row daily figures (Fecha, Calidad,
# Volumen, Precio mínimo/máximo/promedio ponderado, Precio $/Kg) across
# the existing data rows (2-43). Row N after the edit carries the values
# that used to live at row Map[N] before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 31
    3  = 23
    4  = 30
    5  = 5
    6  = 20
    7  = 9
    8  = 14
    9  = 40
    10 = 12
    11 = 41
    12 = 16
    13 = 43
    14 = 25
    15 = 22
    16 = 37
    17 = 21
    18 = 24
    19 = 6
    20 = 17
    21 = 34
    22 = 32
    23 = 11
    24 = 4
    25 = 7
    26 = 8
    27 = 15
    28 = 28
    29 = 19
    30 = 35
    31 = 38
    32 = 42
    33 = 3
    34 = 26
    35 = 27
    36 = 36
    37 = 33
    38 = 10
    39 = 2
    40 = 39
    41 = 18
    42 = 29
    43 = 13
}

# Columns that carry the per-row data that gets reshuffled.
$cols = @("D", "I", "J", "K", "L", "M", "P")

# Snapshot the "before" values for every relevant cell first, since the
# source and destination rows overlap and we must not read values that
# have already been overwritten. Value2 (not Value) is used throughout
# because it returns/accepts plain scalars instead of a Variant wrapper.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 43; $r++) {
        $addr = "$col$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value2 = $snapshot[$srcAddr]
    }
}
